$p = $ppt.ActivePresentation
$full = $p.FullName
$p.Close()

Add-Type -AssemblyName System.IO.Compression.FileSystem

function Get-ZipEntryText($za, $name) {
    $entry = $za.GetEntry($name)
    $stream = $entry.Open()
    $reader = New-Object System.IO.StreamReader($stream, [System.Text.Encoding]::UTF8)
    $content = $reader.ReadToEnd()
    $reader.Close()
    $stream.Close()
    return $content
}

function Set-ZipEntryText($za, $name, $text) {
    $oldEntry = $za.GetEntry($name)
    $oldEntry.Delete()
    $newEntry = $za.CreateEntry($name)
    $stream = $newEntry.Open()
    $writer = New-Object System.IO.StreamWriter($stream, [System.Text.UTF8Encoding]::new($false))
    $writer.Write($text)
    $writer.Flush()
    $writer.Close()
    $stream.Close()
}

$za = [System.IO.Compression.ZipFile]::Open($full, [System.IO.Compression.ZipArchiveMode]::Update)
$t1 = Get-ZipEntryText $za "ppt/theme/theme1.xml"
$t2 = Get-ZipEntryText $za "ppt/theme/theme2.xml"
Set-ZipEntryText $za "ppt/theme/theme1.xml" $t2
Set-ZipEntryText $za "ppt/theme/theme2.xml" $t1
$za.Dispose()

Write-Output "verifying on-disk bytes after dispose"
$za2 = [System.IO.Compression.ZipFile]::Open($full, [System.IO.Compression.ZipArchiveMode]::Read)
$check1 = Get-ZipEntryText $za2 "ppt/theme/theme1.xml"
Write-Output ($check1.Substring(0,150))
$za2.Dispose()

$p2 = $ppt.Presentations.Open($full)
Write-Output ("reopened slides=" + $ppt.ActivePresentation.Slides.Count)
